{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"426\u00f79=47, 3\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"302\u00f75=60, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"881\u00f74=220, 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"207\u00f79=23, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"834\u00f78=104, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"864\u00f75=172, 4\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"636\u00f78=79, 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"144\u00f72=72, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"490\u00f79=54, 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"640\u00f72=320, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"619\u00f77=88, 3\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"303\u00f78=37, 7\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"141\u00f76=23, 3\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"863\u00f77=123, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"963\u00f72=481, 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"397\u00f74=99, 1\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"527\u00f75=105, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"455\u00f76=75, 5\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"676\u00f72=338, 0\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"743\u00f78=92, 7\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"805\u00f74=201, 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"683\u00f77=97, 4\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"900\u00f78=112, 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"106\u00f77=15, 1\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"940\u00f75=188, 0\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"940\u00f77=134, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"328\u00f78=41, 0\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"957\u00f75=191, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"762\u00f72=381, 0\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"799\u00f72=399, 1\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"442\u00f74=110, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"266\u00f79=29, 5\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"289\u00f77=41, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"994\u00f72=497, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"709\u00f73=236, 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"354\u00f75=70, 4\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"162\u00f74=40, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"525\u00f76=87, 3\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"492\u00f77=70, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"587\u00f73=195, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"579\u00f75=115, 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"414\u00f79=46, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"985\u00f78=123, 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"642\u00f74=160, 2\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"183\u00f74=45, 3\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"440\u00f74=110, 0\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"707\u00f75=141, 2\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"694\u00f78=86, 6\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"805\u00f79=89, 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"810\u00f73=270, 0\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"426\u00f79=47, 3\"\n$find.Replacement.Text = \"302\u00f75=60, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"881\u00f74=220, 1\"\n$find.Replacement.Text = \"207\u00f79=23, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"834\u00f78=104, 2\"\n$find.Replacement.Text = \"864\u00f75=172, 4\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"636\u00f78=79, 4\"\n$find.Replacement.Text = \"144\u00f72=72, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"490\u00f79=54, 4\"\n$find.Replacement.Text = \"640\u00f72=320, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"619\u00f77=88, 3\"\n$find.Replacement.Text = \"303\u00f78=37, 7\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"141\u00f76=23, 3\"\n$find.Replacement.Text = \"863\u00f77=123, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"963\u00f72=481, 1\"\n$find.Replacement.Text = \"397\u00f74=99, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"527\u00f75=105, 2\"\n$find.Replacement.Text = \"455\u00f76=75, 5\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"676\u00f72=338, 0\"\n$find.Replacement.Text = \"743\u00f78=92, 7\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"805\u00f74=201, 1\"\n$find.Replacement.Text = \"683\u00f77=97, 4\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"900\u00f78=112, 4\"\n$find.Replacement.Text = \"106\u00f77=15, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"940\u00f75=188, 0\"\n$find.Replacement.Text = \"940\u00f77=134, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"328\u00f78=41, 0\"\n$find.Replacement.Text = \"957\u00f75=191, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"762\u00f72=381, 0\"\n$find.Replacement.Text = \"799\u00f72=399, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"442\u00f74=110, 2\"\n$find.Replacement.Text = \"266\u00f79=29, 5\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"289\u00f77=41, 2\"\n$find.Replacement.Text = \"994\u00f72=497, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"709\u00f73=236, 1\"\n$find.Replacement.Text = \"354\u00f75=70, 4\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"162\u00f74=40, 2\"\n$find.Replacement.Text = \"525\u00f76=87, 3\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"492\u00f77=70, 2\"\n$find.Replacement.Text = \"587\u00f73=195, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"579\u00f75=115, 4\"\n$find.Replacement.Text = \"414\u00f79=46, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"985\u00f78=123, 1\"\n$find.Replacement.Text = \"642\u00f74=160, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"183\u00f74=45, 3\"\n$find.Replacement.Text = \"440\u00f74=110, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"707\u00f75=141, 2\"\n$find.Replacement.Text = \"694\u00f78=86, 6\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"805\u00f79=89, 4\"\n$find.Replacement.Text = \"810\u00f73=270, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n"}
